$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.228090882301331
$ws.Range("B1").Value = 5.544131278991699
$ws.Range("C1").Value = 3.719323396682739
$ws.Range("D1").Value = 0.9812163710594177
$ws.Range("E1").Value = 0.627245306968689
